# "validation added on license page"
#
# The credit-card / post-payment navigation test case (row 4) is replaced
# by a new license-info-entry test case, and the whole suite's execution
# date is bumped from 14-11-2024 to 24-11-2024. The now-redundant last row
# (which duplicated the new license test's data) is removed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 becomes the "verify customer license info entry" test, re-run on
# the newer date.
$ws.Range("A3").Value = "verifyCustomerLicenseInfoEntry"
$ws.Range("C2").Value = "24-11-2024"
$ws.Range("C3").Value = "24-11-2024"

# Row 4 (verifyCustomerNavigationAfterPayment) is no longer needed.
$ws.Rows.Item(4).Delete()

# Column A no longer needs to fit the longest (now-removed) method name.
$ws.Columns.Item(1).AutoFit()
